# Add invoice from latest Arrow order (Arrow Order 5) to expense report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 (Arrow Order 5 / Morgan / 58.73) previously had no receipt file linked.
# Give it the receipt file name, to be turned into a hyperlink below.
$ws.Range("D25").Value = "Arrow Order 5.pdf"

# Describe what was purchased in that Arrow order: PCB manufacturing (row 27)
# and longer programming cables (row 26), in the Items column (F).
$ws.Range("F27").Value = "PCB Manufacturing"
$ws.Range("F26").Value = "Longer Programming Cables"

# Turn the D25 receipt text into a working hyperlink to the PDF, just like
# the other receipt cells in column D.
$null = $ws.Hyperlinks.Add($ws.Range("D25"), "Arrow%20Order%205.pdf", "", "", "Arrow Order 5.pdf")

# Hyperlinks.Add() re-applies formatting; restore the normal Hyperlink cell
# style (matching the other receipt-link cells) instead of any stray style.
$ws.Range("D25").Style = "Hyperlink"

# Reflect the cursor/selection landing on the newly-edited cell, as in the
# saved workbook.
$null = $ws.Range("F26").Select()
